$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("C2").Value = 9521.12606904989
$ws.Range("D2").Value = 16512.26
$ws.Range("F2").Value = -14.8841585110742

# Row 3
$ws.Range("C3").Value = 8936.65268313907
$ws.Range("F3").Value = 192.718558919094

# Row 4
$ws.Range("C4").Value = 6393.20455969602
$ws.Range("F4").Value = 81.8173448307325

# Row 5
$ws.Range("C5").Value = 6627.43916369861
$ws.Range("F5").Value = 99.9361868472916

# Row 6
$ws.Range("C6").Value = 9740.25110897147
$ws.Range("F6").Value = 249.517854039607

# Row 7
$ws.Range("C7").Value = 9493.4729467498
$ws.Range("F7").Value = 9.5998478305631

# Row 9
$ws.Range("C9").Value = 9002.4218949426
$ws.Range("F9").Value = 14.7903675683412

# Row 10
$ws.Range("C10").Value = 8341.97093611705
$ws.Range("F10").Value = 11.2715776172765

# Row 11
$ws.Range("C11").Value = 9172.27400542071
$ws.Range("F11").Value = 13.1666612570396

# Row 12
$ws.Range("C12").Value = 9027.46475388457
$ws.Range("F12").Value = -27.1661799757432

# Row 13
$ws.Range("C13").Value = 9665.69817110578
$ws.Range("F13").Value = -46.1438357486933

# Row 14
$ws.Range("C14").Value = 10130.0475884752
$ws.Range("F14").Value = 304.204056641698

# Row 15
$ws.Range("C15").Value = 10351.002628396
$ws.Range("F15").Value = 313.410516638398
